# Natmi following Dr Hou advice
# Update LR-pairs sheet: recompute cluster-pair rows across ECs / sCs clusters
# (row 2 updated in place; rows 3-7 appended) per the revised NATMI output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 6
$numCols = 20

$data = New-Object 'object[,]' $numRows,$numCols

$data[0,0] = "ECs"
$data[0,1] = "Cd6"
$data[0,2] = "Alcam"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.1475986666666667
$data[0,7] = 0.442796
$data[0,8] = 0.6430044319495352
$data[0,9] = 0.643004431949535
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 62.12558000000001
$data[0,13] = 186.37674
$data[0,14] = 0.9736910227596813
$data[0,15] = 0.9736910227596813
$data[0,16] = 9.169652773893336
$data[0,17] = 82.52687496504001
$data[0,18] = 0.6260876429839508
$data[0,19] = 0.6260876429839507
$data[1,0] = "ECs"
$data[1,1] = "Cd6"
$data[1,2] = "Alcam"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.1475986666666667
$data[1,7] = 0.442796
$data[1,8] = 0.6430044319495352
$data[1,9] = 0.643004431949535
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.5683613333333334
$data[1,13] = 1.705084
$data[1,14] = 0.008907897969731461
$data[1,15] = 0.008907897969731461
$data[1,16] = 0.08388937498488891
$data[1,17] = 0.7550043748640001
$data[1,18] = 0.005727817873891596
$data[1,19] = 0.005727817873891595
$data[2,0] = "ECs"
$data[2,1] = "Cd6"
$data[2,2] = "Alcam"
$data[2,3] = "sCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.1475986666666667
$data[2,7] = 0.442796
$data[2,8] = 0.6430044319495352
$data[2,9] = 0.643004431949535
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 1.110262
$data[2,13] = 3.330786
$data[2,14] = 0.01740107927058724
$data[2,15] = 0.01740107927058724
$data[2,16] = 0.1638731908506667
$data[2,17] = 1.474858717656
$data[2,18] = 0.01118897109169278
$data[2,19] = 0.01118897109169278
$data[3,0] = "sCs"
$data[3,1] = "Cd6"
$data[3,2] = "Alcam"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.08194666666666665
$data[3,7] = 0.24584
$data[3,8] = 0.356995568050465
$data[3,9] = 0.356995568050465
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 62.12558000000001
$data[3,13] = 186.37674
$data[3,14] = 0.9736910227596813
$data[3,15] = 0.9736910227596813
$data[3,16] = 5.090984195733333
$data[3,17] = 45.8188577616
$data[3,18] = 0.3476033797757306
$data[3,19] = 0.3476033797757306
$data[4,0] = "sCs"
$data[4,1] = "Cd6"
$data[4,2] = "Alcam"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.08194666666666665
$data[4,7] = 0.24584
$data[4,8] = 0.356995568050465
$data[4,9] = 0.356995568050465
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.5683613333333334
$data[4,13] = 1.705084
$data[4,14] = 0.008907897969731461
$data[4,15] = 0.008907897969731461
$data[4,16] = 0.04657531672888889
$data[4,17] = 0.41917785056
$data[4,18] = 0.003180080095839866
$data[4,19] = 0.003180080095839866
$data[5,0] = "sCs"
$data[5,1] = "Cd6"
$data[5,2] = "Alcam"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.08194666666666665
$data[5,7] = 0.24584
$data[5,8] = 0.356995568050465
$data[5,9] = 0.356995568050465
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 1.110262
$data[5,13] = 3.330786
$data[5,14] = 0.01740107927058724
$data[5,15] = 0.01740107927058724
$data[5,16] = 0.09098227002666666
$data[5,17] = 0.81884043024
$data[5,18] = 0.006212108178894462
$data[5,19] = 0.006212108178894464

for ($i = 0; $i -lt $numRows; $i++) {
  $targetRow = 2 + $i
  for ($j = 0; $j -lt $numCols; $j++) {
    $targetCol = 1 + $j
    $ws.Cells.Item($targetRow, $targetCol).Value = $data[$i, $j]
  }
}
